$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.947.24'
$ws.Range("E2").Value = '  +5.29%  '
$ws.Range("D3").Value = '2.257.65'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.48%  '
$ws.Range("E7").Value = '  +3.13%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.21%  '
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.114'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").Value = '2.607.35'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").Value = '2.249.01'
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.756'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.00%  '
$ws.Range("D19").Value = '41.826.81'
$ws.Range("E19").Value = '  +5.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.40%  '
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  +3.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.02%  '
$ws.Range("E30").Value = '  +3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0742'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.58%  '
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("E38").Value = '  +4.73%  '
$ws.Range("E39").Value = '  +3.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.00%  '
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("E42").Value = '  +4.49%  '
$ws.Range("D43").Value = '2.049.70'
$ws.Range("E43").Value = '  -3.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.42%  '
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.35%  '
$ws.Range("E49").Value = '  +2.46%  '
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.02%  '
